# Generate Report for Handback
# Updates the Overview status text to reflect a failed handback transform,
# widens the "Error Detail" column on the per-language sheets, and fills
# in the Error Detail message for the affected file on each language sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the Status for the "36f62662..." row (row 3) everywhere it is
# shown: the Overview sheet (zh-cn and de-de columns) and the Status
# column on each language-specific sheet.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Widen the "Error Detail" column (column P) on both language sheets to a
# stored OOXML column width of 40. Excel's ColumnWidth property (character
# units) is offset from the stored width by 5/6, so we compensate here to
# land exactly on width="40" in the saved file.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

# Populate the Error Detail cell for the "36f62662..." row (row 3) on
# each language sheet with the handback/handoff filename mismatch message.
$zhcn.Range("P3").Value = "Handback file name: kmu5droo.j4p is different with handoff file name: 36f62662-67ec-4004-9bae-67d6ac2fe513.c5e74be22a3d910b3dda1e834c5ea4ab418222b8.zh-cn."
$dede.Range("P3").Value = "Handback file name: kmu5droo.j4p is different with handoff file name: 36f62662-67ec-4004-9bae-67d6ac2fe513.c5e74be22a3d910b3dda1e834c5ea4ab418222b8.de-de."
